$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Last status check on" timestamp in F1
$ws.Range("F1").Value = "Last status check on: 10.02.2022 00:30"

# D10: was a text string "+0.3" -> now a real number 0.3
$ws.Range("D10").Value = 0.3

# E10: was a text string "2022-02-10 00:20:52" -> now a real date/time serial
# value, formatted with the same date-time number format used by the other
# rows in column E (e.g. E2:E9).
$ws.Range("E10").Value = 44602.01449074074
$ws.Range("E10").NumberFormat = $ws.Range("E9").NumberFormat
